# Update gh-pages data output (苏州-漫展信息.xlsx)
# Applies the numeric/text refresh captured in the commit "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet index 1)
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 12887
$ws1.Range("F6").Value = 71

$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "2024-09-21"
$ws1.Range("B8").Style = "Normal"
$ws1.Range("E8").Value = "2024.09.21 10:00-09.21 17:00"
$ws1.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202408/PnsN5NWZ1724911969688.jpeg"

$ws1.Range("F9").Value = 13
$ws1.Range("F10").Value = 12813
$ws1.Range("F12").Value = 34
$ws1.Range("F13").Value = 8679
$ws1.Range("F14").Value = 7672
$ws1.Range("F15").Value = 194
$ws1.Range("F16").Value = 106
$ws1.Range("F19").Value = 981
$ws1.Range("F22").Value = 381
$ws1.Range("F25").Value = 86

# Sheet "全部类型" (sheet index 4) - mirrors the same events, offset by one row
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 12887
$ws4.Range("F7").Value = 71

$ws4.Range("B9").NumberFormat = "@"
$ws4.Range("B9").Value = "2024-09-21"
$ws4.Range("B9").Style = "Normal"
$ws4.Range("E9").Value = "2024.09.21 10:00-09.21 17:00"
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202408/PnsN5NWZ1724911969688.jpeg"

$ws4.Range("F10").Value = 13
$ws4.Range("F11").Value = 12813
$ws4.Range("F13").Value = 34
$ws4.Range("F14").Value = 8679
$ws4.Range("F15").Value = 7672
$ws4.Range("F16").Value = 194
$ws4.Range("F17").Value = 106
$ws4.Range("F20").Value = 981
$ws4.Range("F24").Value = 381
$ws4.Range("F27").Value = 86
